$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 165761
$ws.Range("C4").Value = 156700
$ws.Range("C5").Value = 9061
$ws.Range("C8").Value = 65.20999999999999
